# Generate Report for Handoff
#
# The localization-status report was regenerated: the 12d9d1f4-... file
# group (rows 4-7 on every per-locale sheet) got a fresh handoff pass.
#   - Priority moved from "low" to "ht" for that whole batch (zh-cn + de-de).
#   - The handoff timestamp advanced:
#       zh-cn "Latest Handoff Datetime"            04:35:58 -> 04:36:24
#       de-de "Latest Handoff Datetime" /
#         Overview "Latest HO Xliff Generate Date"  04:36:06 -> 04:36:29
#     (de-de and Overview share the same underlying timestamp string, so
#     both locations must be updated together to keep them in sync.)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$rows = 4,5,6,7

foreach ($r in $rows) {
    # Priority: low -> ht
    $ws2.Range("E$r").Value = "ht"
    $ws3.Range("E$r").Value = "ht"

    # zh-cn Latest Handoff Datetime: 2016-08-20 04:35:58 -> 2016-08-20 04:36:24
    $ws2.Range("H$r").Value = "2016-08-20 04:36:24"

    # de-de Latest Handoff Datetime AND the matching Overview "Latest HO
    # Xliff Generate Date" column both held 2016-08-20 04:36:06 -> 04:36:29.
    $ws3.Range("H$r").Value = "2016-08-20 04:36:29"
    $ws1.Range("G$r").Value = "2016-08-20 04:36:29"
}
